$wb = $excel.ActiveWorkbook

# --- Work on the "Collaborative Filtering" sheet (Table2 / table4.xml) ---
$ws = $wb.Worksheets.Item("Collaborative Filtering")

# 1) Mark the two in-progress tasks as finished, with completion dates,
#    by copying the "Good"/date styles already used elsewhere in the table.
$ws.Range("C2").Copy($ws.Range("C6"))
$ws.Range("C6").Value = "Yes"
$ws.Range("D2").Copy($ws.Range("D6"))
$ws.Range("D6").Value = 45319

$ws.Range("C3").Copy($ws.Range("C7"))
$ws.Range("C7").Value = "Yes"
$ws.Range("D3").Copy($ws.Range("D7"))
$ws.Range("D7").Value = 45320

# 2) Insert a new sub-task row above the old row 8, pushing the existing
#    "Create a graph..." row down to row 9 (keeping its formatting intact).
$ws.Rows.Item(8).Insert()

# Give the freshly inserted row 8 the same cell styles used by similar
# "Medium priority / Not finished" rows elsewhere in this table.
$ws.Range("C4").Copy($ws.Range("C8"))
$ws.Range("D4").Copy($ws.Range("D8"))

$ws.Range("A8").Value = "Compare different algorithms e.g. KNNWithMeans, KNNWithZScore etc."
$ws.Range("B8").Value = "Medium"
$ws.Range("C8").Value = "No"

# 3) Append a brand new row 10 for another new sub-task.
$ws.Range("A5").Copy($ws.Range("A10"))
$ws.Range("B5").Copy($ws.Range("B10"))
$ws.Range("C4").Copy($ws.Range("C10"))
$ws.Range("D9").Copy($ws.Range("D10"))

$ws.Range("A10").Value = "Perform a GridSearch to hypertune the algorithm parameters"
$ws.Range("B10").Value = "Low"
$ws.Range("C10").Value = "No"

# 4) Grow the table ("Table2") to cover the two extra rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D10"))

# 5) Make "Collaborative Filtering" the active/selected sheet, matching the
#    saved workbook view (it was the sheet being worked on).
$ws.Select()
$ws.Range("L9").Select()
